$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "91.069.14"
$ws.Range("E2").Value = "  +1.78%  "
$ws.Range("D3").Value = "3.168.20"
$ws.Range("E3").Value = "  +2.31%  "
$ws.Range("E4").Value = "  +0.17%  "
$ws.Range("D5").Value = "'217.34"
$ws.Range("E5").Value = "  +2.06%  "
$ws.Range("D6").Value = "'626.61"
$ws.Range("E6").Value = "  +0.70%  "
$ws.Range("E7").Value = "  +24.85%  "
$ws.Range("D8").Value = "'0.374"
$ws.Range("E8").Value = "  +0.47%  "
$ws.Range("E9").Value = "  -0.02%  "
$ws.Range("D10").Value = "3.158.41"
$ws.Range("E10").Value = "  +2.04%  "
$ws.Range("E11").Value = "  +22.53%  "
$ws.Range("E12").Value = "  +8.21%  "
$ws.Range("E13").Value = "  +3.74%  "
$ws.Range("D14").Value = "'34.97"
$ws.Range("E14").Value = "  +8.47%  "
$ws.Range("E15").Value = "  +4.27%  "
$ws.Range("D16").Value = "90.919.39"
$ws.Range("D17").Value = "3.716.03"
$ws.Range("E17").Value = "  +1.02%  "
$ws.Range("D18").Value = "3.162.05"
$ws.Range("E18").Value = "  +2.00%  "
$ws.Range("E19").Value = "  +11.01%  "
$ws.Range("D20").Value = "'0.0000218"
$ws.Range("E20").Value = "  +1.36%  "
$ws.Range("D21").Value = "'14.30"
$ws.Range("E21").Value = "  +6.13%  "
$ws.Range("D22").Value = "'447.72"
$ws.Range("E22").Value = "  +5.11%  "
$ws.Range("E23").Value = "  +8.03%  "
$ws.Range("E24").Value = "  +5.54%  "
$ws.Range("D25").Value = "'6.14"
$ws.Range("E25").Value = "  +11.28%  "
$ws.Range("D26").Value = "'88.88"
$ws.Range("E26").Value = "  +5.61%  "
$ws.Range("D27").Value = "'12.36"
$ws.Range("E27").Value = "  +2.70%  "
$ws.Range("D28").Value = "3.300.16"
$ws.Range("E28").Value = "  +1.24%  "
$ws.Range("E29").Value = "  -0.06%  "
$ws.Range("E30").Value = "  -0.11%  "
$ws.Range("D31").Value = "'9.16"
$ws.Range("E31").Value = "  +12.50%  "
$ws.Range("D32").Value = "'528.89"
$ws.Range("E32").Value = "  +3.48%  "
$ws.Range("D33").Value = "'0.906"
$ws.Range("E33").Value = "  -14.80%  "
$ws.Range("D34").Value = "'25.81"
$ws.Range("E34").Value = "  +15.00%  "
$ws.Range("D35").Value = "'3.74"
$ws.Range("E35").Value = "  +0.87%  "
$ws.Range("E36").Value = "  +5.23%  "
$ws.Range("E37").Value = "  +10.12%  "
$ws.Range("E38").Value = "  +5.55%  "
$ws.Range("E39").Value = "  +3.99%  "
$ws.Range("D40").Value = "'22.26"
$ws.Range("E40").Value = "  -0.19%  "
$ws.Range("E41").Value = "  -0.09%  "
$ws.Range("E42").Value = "  +15.20%  "
$ws.Range("E43").Value = "  +11.94%  "
$ws.Range("D44").Value = "'0.0821"
$ws.Range("E44").Value = "  +14.77%  "
$ws.Range("D46").Value = "'1.94"
$ws.Range("E46").Value = "  +5.60%  "
$ws.Range("D47").Value = "'148.84"
$ws.Range("E47").Value = "  +2.20%  "
$ws.Range("E48").Value = "  +9.27%  "
$ws.Range("D49").Value = "'44.22"
$ws.Range("E49").Value = "  +1.69%  "
$ws.Range("D50").Value = "'4.43"
$ws.Range("E50").Value = "  +11.64%  "
$ws.Range("D51").Value = "'171.98"
$ws.Range("E51").Value = "  +7.35%  "
